$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SupplierCreation")

# Row 4: PASS -> FAIL
$ws.Range("F4").Value = "FAIL"

# Rows 5-41: clear the stray PASS/FAIL stamps entirely
$ws.Range("F5:F41").ClearContents()
